$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.484.90'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '2.629.74'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '112.78'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.07%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '323.88'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.528'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.15%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.544'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.97'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '19.82'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.99%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0813'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.30'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = '3.038.29'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('D16').Value = '2.626.90'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.861'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.67%  '
$ws.Range('D18').Value = '49.367.65'
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.03'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.95'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.71'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('E22').Value = '  -1.34%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '270.19'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '68.67'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -5.66%  '
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.23'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.44%  '
$ws.Range('E27').Value = '  -1.64%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.34'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.77%  '
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '35.07'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.83%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.138'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.72%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '49.58'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.49'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('E35').Value = '  +2.78%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '19.04'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.62%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.92'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.20%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.05'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.14'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '127.03'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.07%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '22.40'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0322'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.41%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.14'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.05%  '
$ws.Range('D46').Value = '2.061.52'
$ws.Range('E46').Value = '  +0.24%  '
$ws.Range('E47').Value = '  +6.98%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.24'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.90%  '
$ws.Range('E49').Value = '  -8.65%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.94'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.02%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '59.07'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.89%  '
